$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9: "Bancos de imagens" entry ---
$ws.Range("A9").Value = "Bancos de imagens "
$ws.Range("B9").Value = "103 bancos de imagens gratúitos"
$ws.Range("C9").Value = "http://marketingdeconteudo.com/melhores-bancos-de-imagens-gratuitos/"

# --- Row 10: "Comandos GIT" entry ---
$ws.Range("A10").Value = "Comandos GIT"
$ws.Range("B10").Value = "Lista de comandos GIT"
$ws.Range("C10").Value = "https://gist.github.com/leocomelli/2545add34e4fec21ec16"

# Turn the Row 9 link cell into a real hyperlink (matches hyperlinks added to the worksheet)
[void]$ws.Hyperlinks.Add($ws.Range("C9"), "http://marketingdeconteudo.com/melhores-bancos-de-imagens-gratuitos/")

# Re-apply the same cell formatting used by the rest of the table (rows 3-8) to the
# newly used rows as well as the trailing empty styled rows (11-20), re-using the
# existing style entries (border/fill/font) rather than creating new ones.
$ws.Range("A3:C3").Copy()
[void]$ws.Range("A9:C20").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Column widths grew to fit the new, longer content in columns A and C.
$ws.Columns("A").ColumnWidth = 17.7
$ws.Columns("C").ColumnWidth = 68.8

# Leave the cursor on the last-edited cell.
[void]$ws.Range("C10").Select()

Write-Host "Rows 9-20 populated/formatted, hyperlink + column widths updated"
